# Update sexting escalation script lines on the AntonellaJourney sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AntonellaJourney")

$ws.Range("B4").Value  = "let go with me babe... I need you to see this 💜"
$ws.Range("B6").Value  = "I'm cumming... right now... don't look away"
$ws.Range("B7").Value  = "I'm right there babe... don't go anywhere, I need you to watch me finish 💜"
$ws.Range("B8").Value  = "oh my god I can't take it 🖤"
$ws.Range("B9").Value  = "you need to see this... I've never been like this before 💜"
$ws.Range("B11").Value = "I'm about to lose it... you need to see what's happening to me right now"
$ws.Range("B12").Value = "I keep pushing my fingers deeper and moaning into my pillow... god this feels so good 💜"
$ws.Range("B13").Value = "I'm rubbing my clit so hard right now and I can't slow down... my legs are shaking"
$ws.Range("B14").Value = "oh fuck 🖤"
$ws.Range("B15").Value = "look at what you're doing to me... I can't hold back anymore 💜"
$ws.Range("B17").Value = "what would you do if you were here right now babe? I need to hear it"
$ws.Range("B18").Value = "I'm dripping wet and every time I think about you watching me it gets worse 💜"
$ws.Range("B19").Value = "but I literally can't stop touching myself right now... it's like my body won't let me"
$ws.Range("B20").Value = "wow... okay I need a second after that 🖤"
$ws.Range("B21").Value = "I want to show you what you made me feel 💜"
$ws.Range("B23").Value = "I'm lying here and my fingers are starting to wander... I blame you for this babe"
$ws.Range("B24").Value = "my skin is tingling everywhere right now and I can feel my heartbeat getting faster... you're doing something to me"
$ws.Range("B25").Value = "mmm you liked that? that's making me feel way braver than usual 🖤"
